# Update cryptocurrency price/volume table with latest scraped values.
# (Mirrors the GitHub Actions "Updated cryptos list" commit.)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.381.80'
$ws.Range('E2').Value = '  -0.15%  '
$ws.Range('D3').Value = '1.849.17'
$ws.Range('E3').Value = '  -0.01%  '
$ws.Range('D4').Value = "'1.001"
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').Value = "'240.79"
$ws.Range('E5').Value = '  +0.19%  '
$ws.Range('D6').Value = "'0.6309"
$ws.Range('E6').Value = '  +0.47%  '
$ws.Range('D7').Value = "'1.001"
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('D8').Value = "'0.07523"
$ws.Range('E8').Value = '  -2.23%  '
$ws.Range('D9').Value = "'0.2912"
$ws.Range('E9').Value = '  -0.03%  '
$ws.Range('D10').Value = "'24.47"
$ws.Range('E10').Value = '  -1.28%  '
$ws.Range('D11').Value = "'0.07745"
$ws.Range('E11').Value = '  -0.11%  '
$ws.Range('D12').Value = '1.848.35'
$ws.Range('E12').Value = '  +0.14%  '
$ws.Range('D13').Value = "'5.027"
$ws.Range('E13').Value = '  +0.00%  '
$ws.Range('D14').Value = "'0.6816"
$ws.Range('E14').Value = '  +0.05%  '
$ws.Range('D15').Value = "'0.00001040"
$ws.Range('E15').Value = '  -3.78%  '
$ws.Range('D16').Value = "'83.13"
$ws.Range('E16').Value = '  -0.50%  '
$ws.Range('D17').Value = '2.117.33'
$ws.Range('E17').Value = '  -0.03%  '
$ws.Range('D18').Value = "'6.121"
$ws.Range('E18').Value = '  -0.82%  '
$ws.Range('D19').Value = '29.422.40'
$ws.Range('E19').Value = '  -0.06%  '
$ws.Range('D20').Value = "'229.21"
$ws.Range('E20').Value = '  +0.20%  '
$ws.Range('D21').Value = "'12.34"
$ws.Range('E21').Value = '  -0.31%  '
$ws.Range('E22').Value = '  +0.02%  '
$ws.Range('D23').Value = "'7.460"
$ws.Range('E23').Value = '  +0.45%  '
$ws.Range('D24').Value = "'0.9997"
$ws.Range('E24').Value = '  -0.09%  '
$ws.Range('D25').Value = "'159.15"
$ws.Range('E25').Value = '  +1.19%  '
$ws.Range('D26').Value = "'0.1386"
$ws.Range('E26').Value = '  +0.81%  '
$ws.Range('D27').Value = "'8.425"
$ws.Range('E27').Value = '  +0.18%  '
$ws.Range('D28').Value = "'17.59"
$ws.Range('E28').Value = '  -0.73%  '
$ws.Range('E29').Value = '  +5.11%  '
$ws.Range('D30').Value = "'1.477"
$ws.Range('E30').Value = '  +0.92%  '
$ws.Range('D31').Value = "'0.05700"
$ws.Range('E31').Value = '  +0.97%  '
$ws.Range('D32').Value = "'4.137"
$ws.Range('E32').Value = '  +0.41%  '
$ws.Range('D33').Value = "'4.048"
$ws.Range('E33').Value = '  +0.35%  '
$ws.Range('D34').Value = "'1.154"
$ws.Range('E34').Value = '  -0.82%  '
$ws.Range('E35').Value = '  -1.40%  '
$ws.Range('D36').Value = "'0.6953"
$ws.Range('E36').Value = '  -1.32%  '
$ws.Range('D37').Value = "'2.590"
$ws.Range('E37').Value = '  -0.22%  '
$ws.Range('D38').Value = "'2.851"
$ws.Range('E38').Value = '  +3.04%  '
$ws.Range('D39').Value = '1.251.47'
$ws.Range('E39').Value = '  +1.86%  '
$ws.Range('D40').Value = "'0.01828"
$ws.Range('E40').Value = '  +2.24%  '
$ws.Range('D41').Value = "'6.498"
$ws.Range('E41').Value = '  +0.67%  '
$ws.Range('D42').Value = "'0.9064"
$ws.Range('E42').Value = '  +0.10%  '
$ws.Range('D43').Value = "'1.000"
$ws.Range('E43').Value = '  -0.12%  '
$ws.Range('D44').Value = '2.014.13'
$ws.Range('E44').Value = '  -0.48%  '
$ws.Range('D45').Value = "'101.48"
$ws.Range('E45').Value = '  -0.35%  '
$ws.Range('D46').Value = "'65.97"
$ws.Range('E46').Value = '  +0.14%  '
$ws.Range('D47').Value = "'7.088"
$ws.Range('E47').Value = '  -1.14%  '
$ws.Range('D48').Value = "'0.1166"
$ws.Range('E48').Value = '  +0.77%  '
$ws.Range('D49').Value = "'9.030"
$ws.Range('E49').Value = '  +0.08%  '
$ws.Range('B50').Value = 'TheSandbox'
$ws.Range('C50').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D50').Value = "'0.3958"
$ws.Range('E50').Value = '  -1.39%  '
$ws.Range('B51').Value = 'RenderToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D51').Value = "'1.669"
$ws.Range('E51').Value = '  -0.41%  '
